# PEMENDEKAN-Lema.xlsx edit
# 1. Row 39, column D: "pa,pak" -> "yah, pa, pak"
# 2. Entire row 51 (Lema "yah" / Sinonim "yah, pa, pak") is deleted, shifting
#    every row below it up by one (so the old row 449 disappears and the
#    sheet now ends at row 448).
# 3. Restore the saved selection/active cell to match the post-edit view
#    (user had scrolled to row 34 and selected F35 before saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D39").Value = "yah, pa, pak"

$ws.Rows(51).Delete()

$ws.Range("F35").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 34 | Out-Null
